$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 157.2
$ws.Range("I9").Value = 149
$ws.Range("J9").Value = 169.5
$ws.Range("K9").Value = 149
$ws.Range("L9").Value = 169.5
$ws.Range("M9").Value = 20
$ws.Range("N9").Value = -507.5
$ws.Range("H17").Value = 2392.25
$ws.Range("I17").Value = 2874.6667
$ws.Range("J17").Value = 1997.5454
$ws.Range("K17").Value = 8624.000100000001
$ws.Range("L17").Value = 5992.6362
$ws.Range("M17").Value = -8456.000100000001
$ws.Range("N17").Value = -6328.6362
$ws.Range("H55").Value = 658.75
$ws.Range("I55").Value = 150
$ws.Range("J55").Value = 731.4286
$ws.Range("K55").Value = 150
$ws.Range("L55").Value = 731.4286
$ws.Range("M55").Value = 64
$ws.Range("N55").Value = -1159.4286
$ws.Range("H131").Value = 5080.0713
$ws.Range("I131").Value = 854.1667
$ws.Range("J131").Value = 8249.5
$ws.Range("K131").Value = 2562.5001
$ws.Range("L131").Value = 24748.5
$ws.Range("M131").Value = 2477.4999
$ws.Range("N131").Value = -34828.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 374
$ws.Range("I2").Value = 186.33333
$ws.Range("J2").Value = 1500
$ws.Range("K2").Value = 186.33333
$ws.Range("L2").Value = 1500
$ws.Range("M2").Value = -73.33332999999999
$ws.Range("N2").Value = -1726
$ws.Range("M17").Value = 165.3333335
$ws.Range("H32").Value = 15941
$ws.Range("I32").Value = 13998.5
$ws.Range("J32").Value = 21380
$ws.Range("K32").Value = 13998.5
$ws.Range("L32").Value = 21380
$ws.Range("M32").Value = -13711.5
$ws.Range("N32").Value = -21954
$ws.Range("H116").Value = 374
$ws.Range("I116").Value = 186.33333
$ws.Range("J116").Value = 1500
$ws.Range("K116").Value = 186.33333
$ws.Range("L116").Value = 1500
$ws.Range("M116").Value = 2107.66667
$ws.Range("N116").Value = -6088
$ws.Range("H132").Value = 3206.32
$ws.Range("I132").Value = 2579.2856
$ws.Range("J132").Value = 6498.25
$ws.Range("K132").Value = 7737.8568
$ws.Range("L132").Value = 19494.75
$ws.Range("M132").Value = -5207.8568
$ws.Range("N132").Value = -24554.75
$ws.Range("N17").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 374
$ws.Range("I3").Value = 186.33333
$ws.Range("J3").Value = 1500
$ws.Range("K3").Value = 186.33333
$ws.Range("L3").Value = 1500
$ws.Range("M3").Value = -72.33332999999999
$ws.Range("N3").Value = -1728
$ws.Range("H10").Value = 5
$ws.Range("I10").Value = 5
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 5
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = 135
$ws.Range("H134").Value = 2446.1765
$ws.Range("I134").Value = 2184.7144
$ws.Range("J134").Value = 3666.3333
$ws.Range("K134").Value = 6554.1432
$ws.Range("L134").Value = 10998.9999
$ws.Range("M134").Value = -4019.1432
$ws.Range("N134").Value = -16068.9999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4694.75
$ws.Range("I31").Value = 2063.7778
$ws.Range("J31").Value = 8077.4287
$ws.Range("K31").Value = 2063.7778
$ws.Range("L31").Value = 8077.4287
$ws.Range("M31").Value = -1768.7778
$ws.Range("N31").Value = -8667.4287
$ws.Range("H34").Value = 4694.75
$ws.Range("I34").Value = 2063.7778
$ws.Range("J34").Value = 8077.4287
$ws.Range("K34").Value = 2063.7778
$ws.Range("L34").Value = 8077.4287
$ws.Range("M34").Value = -1861.7778
$ws.Range("N34").Value = -8481.4287
$ws.Range("H58").Value = 3674.0833
$ws.Range("I58").Value = 1914.75
$ws.Range("J58").Value = 7192.75
$ws.Range("K58").Value = 1914.75
$ws.Range("L58").Value = 7192.75
$ws.Range("M58").Value = -1711.75
$ws.Range("N58").Value = -7598.75
$ws.Range("H62").Value = 135332.67
$ws.Range("I62").Value = 2999
$ws.Range("J62").Value = 201499.5
$ws.Range("K62").Value = 2999
$ws.Range("L62").Value = 201499.5
$ws.Range("M62").Value = -2375
$ws.Range("N62").Value = -202747.5
$ws.Range("H65").Value = 135332.67
$ws.Range("I65").Value = 2999
$ws.Range("J65").Value = 201499.5
$ws.Range("K65").Value = 14995
$ws.Range("L65").Value = 1007497.5
$ws.Range("M65").Value = -11875
$ws.Range("N65").Value = -1013737.5
$ws.Range("H99").Value = 13789.546
$ws.Range("I99").Value = 11277.375
$ws.Range("J99").Value = 15225.071
$ws.Range("K99").Value = 11277.375
$ws.Range("L99").Value = 15225.071
$ws.Range("M99").Value = -9779.375
$ws.Range("N99").Value = -18221.071
$ws.Range("H126").Value = 13789.546
$ws.Range("I126").Value = 11277.375
$ws.Range("J126").Value = 15225.071
$ws.Range("K126").Value = 33832.125
$ws.Range("L126").Value = 45675.213
$ws.Range("M126").Value = -31362.125
$ws.Range("N126").Value = -50615.213
$ws.Range("H132").Value = 2440.6155
$ws.Range("I132").Value = 1617.3636
$ws.Range("J132").Value = 6968.5
$ws.Range("K132").Value = 4852.0908
$ws.Range("L132").Value = 20905.5
$ws.Range("M132").Value = -2322.0908
$ws.Range("N132").Value = -25965.5
$ws.Range("H134").Value = 3924.4285
$ws.Range("I134").Value = 3197.1428
$ws.Range("J134").Value = 4651.7144
$ws.Range("K134").Value = 9591.428400000001
$ws.Range("L134").Value = 13955.1432
$ws.Range("M134").Value = -7056.428400000001
$ws.Range("N134").Value = -19025.1432
$ws.Range("H136").Value = 3674.0833
$ws.Range("I136").Value = 1914.75
$ws.Range("J136").Value = 7192.75
$ws.Range("K136").Value = 5744.25
$ws.Range("L136").Value = 21578.25
$ws.Range("M136").Value = -3194.25
$ws.Range("N136").Value = -26678.25

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 37.5
$ws.Range("I7").Value = 49.333332
$ws.Range("J7").Value = 25.666666
$ws.Range("K7").Value = 147.999996
$ws.Range("L7").Value = 76.99999800000001
$ws.Range("M7").Value = -35.99999600000001
$ws.Range("N7").Value = -300.999998
$ws.Range("H14").Value = 871.2727
$ws.Range("I14").Value = 871.2727
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 2613.8181
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -2440.8181
$ws.Range("L52").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("H75").Value = 359.5
$ws.Range("I75").Value = 388.4
$ws.Range("J75").Value = 215
$ws.Range("K75").Value = 1165.2
$ws.Range("L75").Value = 645
$ws.Range("M75").Value = -167.1999999999998
$ws.Range("N75").Value = -2641
$ws.Range("H78").Value = 359.5
$ws.Range("I78").Value = 388.4
$ws.Range("J78").Value = 215
$ws.Range("K78").Value = 3495.6
$ws.Range("L78").Value = 1935
$ws.Range("M78").Value = 1496.4
$ws.Range("N78").Value = -11919
$ws.Range("L82").Value = 60000
$ws.Range("N82").Value = -60812
$ws.Range("L85").Value = 60000
$ws.Range("N85").Value = -62808
$ws.Range("N52").ClearContents()
$ws.Range("N70").ClearContents()
$ws.Range("N73").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 4756.5557
$ws.Range("I113").Value = 4005
$ws.Range("J113").Value = 4971.2856
$ws.Range("K113").Value = 4005
$ws.Range("L113").Value = 4971.2856
$ws.Range("M113").Value = -1835
$ws.Range("N113").Value = -9311.285599999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4680.273
$ws.Range("I22").Value = 4897.8
$ws.Range("J22").Value = 4499
$ws.Range("K22").Value = 4897.8
$ws.Range("L22").Value = 4499
$ws.Range("M22").Value = -4602.8
$ws.Range("N22").Value = -5089
$ws.Range("H27").Value = 4680.273
$ws.Range("I27").Value = 4897.8
$ws.Range("J27").Value = 4499
$ws.Range("K27").Value = 4897.8
$ws.Range("L27").Value = 4499
$ws.Range("M27").Value = -4790.8
$ws.Range("N27").Value = -4713
$ws.Range("H46").Value = 3604.0527
$ws.Range("I46").Value = 2498
$ws.Range("J46").Value = 3999.0715
$ws.Range("K46").Value = 2498
$ws.Range("L46").Value = 3999.0715
$ws.Range("M46").Value = -2310
$ws.Range("N46").Value = -4375.0715
$ws.Range("H132").Value = 3843.5217
$ws.Range("I132").Value = 3645.05
$ws.Range("J132").Value = 5166.6665
$ws.Range("K132").Value = 10935.15
$ws.Range("L132").Value = 15499.9995
$ws.Range("M132").Value = -8405.150000000001
$ws.Range("N132").Value = -20559.9995

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("L39").Value = 0
$ws.Range("H81").Value = 873.1
$ws.Range("I81").Value = 931.2222
$ws.Range("J81").Value = 350
$ws.Range("K81").Value = 1862.4444
$ws.Range("L81").Value = 700
$ws.Range("M81").Value = -801.4444000000001
$ws.Range("N81").Value = -2822
$ws.Range("H84").Value = 873.1
$ws.Range("I84").Value = 931.2222
$ws.Range("J84").Value = 350
$ws.Range("K84").Value = 9312.222
$ws.Range("L84").Value = 3500
$ws.Range("M84").Value = -4008.222
$ws.Range("N84").Value = -14108
$ws.Range("H107").Value = 1335.2858
$ws.Range("I107").Value = 450
$ws.Range("J107").Value = 1482.8334
$ws.Range("K107").Value = 1350
$ws.Range("L107").Value = 4448.5002
$ws.Range("M107").Value = 570
$ws.Range("N107").Value = -8288.5002
$ws.Range("H113").Value = 935.5333000000001
$ws.Range("I113").Value = 814.9
$ws.Range("J113").Value = 1176.8
$ws.Range("K113").Value = 2444.7
$ws.Range("L113").Value = 3530.4
$ws.Range("M113").Value = -274.6999999999998
$ws.Range("N113").Value = -7870.4
$ws.Range("H132").Value = 1371.1428
$ws.Range("I132").Value = 1092
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 3276
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -746
$ws.Range("N132").Value = -20060
$ws.Range("H136").Value = 2908.6191
$ws.Range("I136").Value = 1102.6154
$ws.Range("J136").Value = 5843.375
$ws.Range("K136").Value = 2908.6191
$ws.Range("L136").Value = 17530.125
$ws.Range("M136").Value = -757.8462
$ws.Range("N136").Value = -22630.125
$ws.Range("M39").ClearContents()
$ws.Range("N39").ClearContents()
